$d = $word.ActiveDocument

$d.Content.Find.Execute("23-5=18", $true, $false, $false, $false, $false, $true, 1, $false, "35-24=11", 2) | Out-Null
$d.Content.Find.Execute("41-14=27", $true, $false, $false, $false, $false, $true, 1, $false, "44+0=44", 2) | Out-Null
$d.Content.Find.Execute("53-15=38", $true, $false, $false, $false, $false, $true, 1, $false, "12+73=85", 2) | Out-Null
$d.Content.Find.Execute("0+44=44", $true, $false, $false, $false, $false, $true, 1, $false, "49-10=39", 2) | Out-Null
$d.Content.Find.Execute("86+5=91", $true, $false, $false, $false, $false, $true, 1, $false, "46-5=41", 2) | Out-Null
$d.Content.Find.Execute("86-57=29", $true, $false, $false, $false, $false, $true, 1, $false, "70-7=63", 2) | Out-Null
$d.Content.Find.Execute("36-34=2", $true, $false, $false, $false, $false, $true, 1, $false, "54+13=67", 2) | Out-Null
$d.Content.Find.Execute("82-75=7", $true, $false, $false, $false, $false, $true, 1, $false, "32-13=19", 2) | Out-Null
$d.Content.Find.Execute("39+35=74", $true, $false, $false, $false, $false, $true, 1, $false, "62+26=88", 2) | Out-Null
$d.Content.Find.Execute("56+41=97", $true, $false, $false, $false, $false, $true, 1, $false, "21+22=43", 2) | Out-Null
$d.Content.Find.Execute("72-21=51", $true, $false, $false, $false, $false, $true, 1, $false, "78-71=7", 2) | Out-Null
$d.Content.Find.Execute("64-30=34", $true, $false, $false, $false, $false, $true, 1, $false, "82-76=6", 2) | Out-Null
$d.Content.Find.Execute("7+43=50", $true, $false, $false, $false, $false, $true, 1, $false, "20-0=20", 2) | Out-Null
$d.Content.Find.Execute("14+8=22", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=6", 2) | Out-Null
$d.Content.Find.Execute("51-16=35", $true, $false, $false, $false, $false, $true, 1, $false, "76-30=46", 2) | Out-Null
$d.Content.Find.Execute("95-81=14", $true, $false, $false, $false, $false, $true, 1, $false, "54+10=64", 2) | Out-Null
$d.Content.Find.Execute("9+24=33", $true, $false, $false, $false, $false, $true, 1, $false, "14+76=90", 2) | Out-Null
$d.Content.Find.Execute("6-1=5", $true, $false, $false, $false, $false, $true, 1, $false, "8+30=38", 2) | Out-Null
$d.Content.Find.Execute("40-14=26", $true, $false, $false, $false, $false, $true, 1, $false, "71-63=8", 2) | Out-Null
$d.Content.Find.Execute("35+30=65", $true, $false, $false, $false, $false, $true, 1, $false, "17+35=52", 2) | Out-Null
$d.Content.Find.Execute("28+38=66", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=7", 2) | Out-Null
$d.Content.Find.Execute("39+33=72", $true, $false, $false, $false, $false, $true, 1, $false, "58-50=8", 2) | Out-Null
$d.Content.Find.Execute("46+23=69", $true, $false, $false, $false, $false, $true, 1, $false, "29-27=2", 2) | Out-Null
$d.Content.Find.Execute("35-23=12", $true, $false, $false, $false, $false, $true, 1, $false, "14+13=27", 2) | Out-Null
$d.Content.Find.Execute("63+11=74", $true, $false, $false, $false, $false, $true, 1, $false, "38-25=13", 2) | Out-Null
$d.Content.Find.Execute("19+27=46", $true, $false, $false, $false, $false, $true, 1, $false, "8+18=26", 2) | Out-Null
$d.Content.Find.Execute("61-34=27", $true, $false, $false, $false, $false, $true, 1, $false, "95-60=35", 2) | Out-Null
$d.Content.Find.Execute("35-1=34", $true, $false, $false, $false, $false, $true, 1, $false, "81-22=59", 2) | Out-Null
$d.Content.Find.Execute("98-3=95", $true, $false, $false, $false, $false, $true, 1, $false, "3+83=86", 2) | Out-Null
$d.Content.Find.Execute("92-10=82", $true, $false, $false, $false, $false, $true, 1, $false, "51-14=37", 2) | Out-Null
$d.Content.Find.Execute("77-21=56", $true, $false, $false, $false, $false, $true, 1, $false, "2+73=75", 2) | Out-Null
$d.Content.Find.Execute("5+14=19", $true, $false, $false, $false, $false, $true, 1, $false, "82-33=49", 2) | Out-Null
$d.Content.Find.Execute("15+70=85", $true, $false, $false, $false, $false, $true, 1, $false, "79-72=7", 2) | Out-Null
$d.Content.Find.Execute("96-68=28", $true, $false, $false, $false, $false, $true, 1, $false, "53+46=99", 2) | Out-Null
$d.Content.Find.Execute("71-60=11", $true, $false, $false, $false, $false, $true, 1, $false, "35+38=73", 2) | Out-Null
$d.Content.Find.Execute("27-11=16", $true, $false, $false, $false, $false, $true, 1, $false, "37-15=22", 2) | Out-Null
$d.Content.Find.Execute("23+26=49", $true, $false, $false, $false, $false, $true, 1, $false, "92-29=63", 2) | Out-Null
$d.Content.Find.Execute("94-68=26", $true, $false, $false, $false, $false, $true, 1, $false, "79-58=21", 2) | Out-Null
$d.Content.Find.Execute("26+22=48", $true, $false, $false, $false, $false, $true, 1, $false, "61-57=4", 2) | Out-Null
$d.Content.Find.Execute("42+13=55", $true, $false, $false, $false, $false, $true, 1, $false, "76-65=11", 2) | Out-Null
$d.Content.Find.Execute("42-4=38", $true, $false, $false, $false, $false, $true, 1, $false, "54-38=16", 2) | Out-Null
$d.Content.Find.Execute("80-50=30", $true, $false, $false, $false, $false, $true, 1, $false, "55-27=28", 2) | Out-Null
$d.Content.Find.Execute("90-48=42", $true, $false, $false, $false, $false, $true, 1, $false, "5+54=59", 2) | Out-Null
$d.Content.Find.Execute("38-21=17", $true, $false, $false, $false, $false, $true, 1, $false, "89-72=17", 2) | Out-Null
$d.Content.Find.Execute("0+38=38", $true, $false, $false, $false, $false, $true, 1, $false, "14+43=57", 2) | Out-Null
$d.Content.Find.Execute("63-2=61", $true, $false, $false, $false, $false, $true, 1, $false, "22-6=16", 2) | Out-Null
$d.Content.Find.Execute("62-32=30", $true, $false, $false, $false, $false, $true, 1, $false, "62-39=23", 2) | Out-Null
$d.Content.Find.Execute("64-21=43", $true, $false, $false, $false, $false, $true, 1, $false, "63+27=90", 2) | Out-Null
$d.Content.Find.Execute("35+2=37", $true, $false, $false, $false, $false, $true, 1, $false, "8+48=56", 2) | Out-Null
$d.Content.Find.Execute("82-69=13", $true, $false, $false, $false, $false, $true, 1, $false, "10+63=73", 2) | Out-Null
$d.Content.Find.Execute("60-11=49", $true, $false, $false, $false, $false, $true, 1, $false, "40-20=20", 2) | Out-Null
$d.Content.Find.Execute("82-18=64", $true, $false, $false, $false, $false, $true, 1, $false, "57+5=62", 2) | Out-Null
$d.Content.Find.Execute("37+9=46", $true, $false, $false, $false, $false, $true, 1, $false, "53+30=83", 2) | Out-Null
$d.Content.Find.Execute("9+82=91", $true, $false, $false, $false, $false, $true, 1, $false, "84-5=79", 2) | Out-Null
$d.Content.Find.Execute("22-1=21", $true, $false, $false, $false, $false, $true, 1, $false, "99-15=84", 2) | Out-Null
$d.Content.Find.Execute("75-17=58", $true, $false, $false, $false, $false, $true, 1, $false, "88-71=17", 2) | Out-Null
$d.Content.Find.Execute("13+30=43", $true, $false, $false, $false, $false, $true, 1, $false, "52-36=16", 2) | Out-Null
$d.Content.Find.Execute("20+70=90", $true, $false, $false, $false, $false, $true, 1, $false, "80-25=55", 2) | Out-Null
$d.Content.Find.Execute("22+39=61", $true, $false, $false, $false, $false, $true, 1, $false, "18-0=18", 2) | Out-Null
$d.Content.Find.Execute("64+10=74", $true, $false, $false, $false, $false, $true, 1, $false, "79-51=28", 2) | Out-Null
$d.Content.Find.Execute("31-31=0", $true, $false, $false, $false, $false, $true, 1, $false, "90-49=41", 2) | Out-Null
$d.Content.Find.Execute("36-15=21", $true, $false, $false, $false, $false, $true, 1, $false, "84-20=64", 2) | Out-Null
$d.Content.Find.Execute("42+7=49", $true, $false, $false, $false, $false, $true, 1, $false, "60-45=15", 2) | Out-Null
$d.Content.Find.Execute("62-33=29", $true, $false, $false, $false, $false, $true, 1, $false, "54-3=51", 2) | Out-Null
$d.Content.Find.Execute("63-39=24", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=25", 2) | Out-Null
$d.Content.Find.Execute("62-36=26", $true, $false, $false, $false, $false, $true, 1, $false, "38+4=42", 2) | Out-Null
$d.Content.Find.Execute("36+34=70", $true, $false, $false, $false, $false, $true, 1, $false, "3+62=65", 2) | Out-Null
$d.Content.Find.Execute("49-22=27", $true, $false, $false, $false, $false, $true, 1, $false, "76-9=67", 2) | Out-Null
$d.Content.Find.Execute("50-44=6", $true, $false, $false, $false, $false, $true, 1, $false, "56-11=45", 2) | Out-Null
$d.Content.Find.Execute("57+13=70", $true, $false, $false, $false, $false, $true, 1, $false, "4+6=10", 2) | Out-Null
$d.Content.Find.Execute("52-51=1", $true, $false, $false, $false, $false, $true, 1, $false, "10+16=26", 2) | Out-Null
$d.Content.Find.Execute("69-67=2", $true, $false, $false, $false, $false, $true, 1, $false, "99-46=53", 2) | Out-Null
$d.Content.Find.Execute("2+10=12", $true, $false, $false, $false, $false, $true, 1, $false, "21+75=96", 2) | Out-Null
$d.Content.Find.Execute("5+68=73", $true, $false, $false, $false, $false, $true, 1, $false, "4-1=3", 2) | Out-Null
$d.Content.Find.Execute("62-3=59", $true, $false, $false, $false, $false, $true, 1, $false, "49-35=14", 2) | Out-Null
$d.Content.Find.Execute("98-89=9", $true, $false, $false, $false, $false, $true, 1, $false, "32+42=74", 2) | Out-Null
$d.Content.Find.Execute("62-22=40", $true, $false, $false, $false, $false, $true, 1, $false, "22+30=52", 2) | Out-Null
$d.Content.Find.Execute("94-10=84", $true, $false, $false, $false, $false, $true, 1, $false, "20+71=91", 2) | Out-Null
$d.Content.Find.Execute("38+13=51", $true, $false, $false, $false, $false, $true, 1, $false, "43-40=3", 2) | Out-Null
$d.Content.Find.Execute("50-40=10", $true, $false, $false, $false, $false, $true, 1, $false, "21+21=42", 2) | Out-Null
$d.Content.Find.Execute("98-78=20", $true, $false, $false, $false, $false, $true, 1, $false, "89-43=46", 2) | Out-Null
$d.Content.Find.Execute("51+1=52", $true, $false, $false, $false, $false, $true, 1, $false, "41+48=89", 2) | Out-Null
$d.Content.Find.Execute("20+56=76", $true, $false, $false, $false, $false, $true, 1, $false, "21+8=29", 2) | Out-Null
$d.Content.Find.Execute("36+52=88", $true, $false, $false, $false, $false, $true, 1, $false, "95-26=69", 2) | Out-Null
$d.Content.Find.Execute("72-57=15", $true, $false, $false, $false, $false, $true, 1, $false, "12+39=51", 2) | Out-Null
$d.Content.Find.Execute("71-37=34", $true, $false, $false, $false, $false, $true, 1, $false, "88-47=41", 2) | Out-Null
$d.Content.Find.Execute("34+35=69", $true, $false, $false, $false, $false, $true, 1, $false, "5+52=57", 2) | Out-Null
$d.Content.Find.Execute("47-45=2", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=36", 2) | Out-Null
$d.Content.Find.Execute("63+34=97", $true, $false, $false, $false, $false, $true, 1, $false, "89-8=81", 2) | Out-Null
$d.Content.Find.Execute("60+0=60", $true, $false, $false, $false, $false, $true, 1, $false, "54+35=89", 2) | Out-Null
$d.Content.Find.Execute("95-61=34", $true, $false, $false, $false, $false, $true, 1, $false, "40-26=14", 2) | Out-Null
$d.Content.Find.Execute("74+6=80", $true, $false, $false, $false, $false, $true, 1, $false, "33-14=19", 2) | Out-Null
$d.Content.Find.Execute("97-94=3", $true, $false, $false, $false, $false, $true, 1, $false, "25+18=43", 2) | Out-Null
$d.Content.Find.Execute("8+50=58", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=55", 2) | Out-Null
$d.Content.Find.Execute("50-47=3", $true, $false, $false, $false, $false, $true, 1, $false, "76-59=17", 2) | Out-Null
$d.Content.Find.Execute("67-19=48", $true, $false, $false, $false, $false, $true, 1, $false, "75-66=9", 2) | Out-Null
$d.Content.Find.Execute("41+47=88", $true, $false, $false, $false, $false, $true, 1, $false, "87-62=25", 2) | Out-Null
$d.Content.Find.Execute("96-15=81", $true, $false, $false, $false, $false, $true, 1, $false, "63+36=99", 2) | Out-Null
$d.Content.Find.Execute("10+23=33", $true, $false, $false, $false, $false, $true, 1, $false, "56-42=14", 2) | Out-Null
$d.Content.Find.Execute("0+94=94", $true, $false, $false, $false, $false, $true, 1, $false, "70+26=96", 2) | Out-Null
